$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3462.6086
$ws.Range("I64").Value = 3210.2703
$ws.Range("J64").Value = 4500
$ws.Range("K64").Value = 3210.2703
$ws.Range("L64").Value = 4500
$ws.Range("M64").Value = -2962.2703
$ws.Range("N64").Value = -4996
# Row 67
$ws.Range("H67").Value = 3462.6086
$ws.Range("I67").Value = 3210.2703
$ws.Range("J67").Value = 4500
$ws.Range("K67").Value = 3210.2703
$ws.Range("L67").Value = 4500
$ws.Range("M67").Value = -2352.2703
$ws.Range("N67").Value = -6216
# Row 107
$ws.Range("H107").Value = 562373.9
$ws.Range("I107").Value = 919293.0600000001
$ws.Range("J107").Value = 1500.8572
$ws.Range("K107").Value = 919293.0600000001
$ws.Range("L107").Value = 1500.8572
$ws.Range("M107").Value = -917373.0600000001
$ws.Range("N107").Value = -5340.8572
# Row 113
$ws.Range("H113").Value = 6000.3945
$ws.Range("I113").Value = 2645.8635
$ws.Range("J113").Value = 10612.875
$ws.Range("K113").Value = 2645.8635
$ws.Range("L113").Value = 10612.875
$ws.Range("M113").Value = 608.1365000000001
$ws.Range("N113").Value = -17120.875
# Row 132
$ws.Range("H132").Value = 28890.447
$ws.Range("I132").Value = 28890.447
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 86671.341
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -84141.341
$ws.Range("N132").ClearContents()
# Row 137
$ws.Range("H137").Value = 653.96
$ws.Range("I137").Value = 663.9091
$ws.Range("J137").Value = 581
$ws.Range("K137").Value = 1991.7273
$ws.Range("L137").Value = 1743
$ws.Range("M137").Value = 558.2727
$ws.Range("N137").Value = -6843

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 874.8444
$ws.Range("I2").Value = 901.4865
$ws.Range("J2").Value = 751.625
$ws.Range("K2").Value = 901.4865
$ws.Range("L2").Value = 751.625
$ws.Range("M2").Value = -788.4865
$ws.Range("N2").Value = -977.625
# Row 32
$ws.Range("H32").Value = 6276.946
$ws.Range("I32").Value = 5177.537
$ws.Range("J32").Value = 16799.857
$ws.Range("K32").Value = 5177.537
$ws.Range("L32").Value = 16799.857
$ws.Range("M32").Value = -4890.537
$ws.Range("N32").Value = -17373.857
# Row 45
$ws.Range("H45").Value = 1937.4286
$ws.Range("I45").Value = 1890.5
$ws.Range("K45").Value = 1890.5
$ws.Range("M45").Value = -1513.5
# Row 61
$ws.Range("H61").Value = 2603.7222
$ws.Range("I61").Value = 2488.9167
$ws.Range("J61").Value = 2833.3333
$ws.Range("K61").Value = 2488.9167
$ws.Range("L61").Value = 2833.3333
$ws.Range("M61").Value = -2276.9167
$ws.Range("N61").Value = -3257.3333
# Row 74
$ws.Range("H74").Value = 5753.0713
$ws.Range("I74").Value = 8651.467000000001
$ws.Range("J74").Value = 2408.7693
$ws.Range("K74").Value = 8651.467000000001
$ws.Range("L74").Value = 2408.7693
$ws.Range("M74").Value = -7777.467000000001
$ws.Range("N74").Value = -4156.7693
# Row 77
$ws.Range("H77").Value = 5753.0713
$ws.Range("I77").Value = 8651.467000000001
$ws.Range("J77").Value = 2408.7693
$ws.Range("K77").Value = 43257.33500000001
$ws.Range("L77").Value = 12043.8465
$ws.Range("M77").Value = -38889.33500000001
$ws.Range("N77").Value = -20779.8465
# Row 97
$ws.Range("H97").Value = 704.61536
$ws.Range("I97").Value = 622.5
$ws.Range("J97").Value = 836
$ws.Range("K97").Value = 622.5
$ws.Range("L97").Value = 836
$ws.Range("M97").Value = -126.5
$ws.Range("N97").Value = -1828
# Row 110
$ws.Range("H110").Value = 56340.11
$ws.Range("I110").Value = 56340.11
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 56340.11
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -54295.11
$ws.Range("N110").ClearContents()
# Row 116
$ws.Range("H116").Value = 874.8444
$ws.Range("I116").Value = 901.4865
$ws.Range("J116").Value = 751.625
$ws.Range("K116").Value = 901.4865
$ws.Range("L116").Value = 751.625
$ws.Range("M116").Value = 1392.5135
$ws.Range("N116").Value = -5339.625
# Row 136
$ws.Range("H136").Value = 2603.7222
$ws.Range("I136").Value = 2488.9167
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 7466.750100000001
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -4916.750100000001
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 874.8444
$ws.Range("I3").Value = 901.4865
$ws.Range("J3").Value = 751.625
$ws.Range("K3").Value = 901.4865
$ws.Range("L3").Value = 751.625
$ws.Range("M3").Value = -787.4865
$ws.Range("N3").Value = -979.625
# Row 94
$ws.Range("H94").Value = 684.1429000000001
$ws.Range("I94").Value = 836.3333
$ws.Range("J94").Value = 570
$ws.Range("K94").Value = 836.3333
$ws.Range("L94").Value = 570
$ws.Range("M94").Value = -385.3333
$ws.Range("N94").Value = -1472
# Row 130
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
# Row 132
$ws.Range("H132").Value = 48772.5
$ws.Range("J132").Value = 48772.5
$ws.Range("L132").Value = 48772.5
$ws.Range("N132").Value = -58892.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1775.5979
$ws.Range("I31").Value = 996.75757
$ws.Range("J31").Value = 3433.7742
$ws.Range("K31").Value = 996.75757
$ws.Range("L31").Value = 3433.7742
$ws.Range("M31").Value = -701.75757
$ws.Range("N31").Value = -4023.7742
# Row 34
$ws.Range("H34").Value = 1775.5979
$ws.Range("I34").Value = 996.75757
$ws.Range("J34").Value = 3433.7742
$ws.Range("K34").Value = 996.75757
$ws.Range("L34").Value = 3433.7742
$ws.Range("M34").Value = -794.75757
$ws.Range("N34").Value = -3837.7742
# Row 58
$ws.Range("H58").Value = 2417.25
$ws.Range("I58").Value = 2212.8
$ws.Range("K58").Value = 2212.8
$ws.Range("M58").Value = -2009.8
# Row 94
$ws.Range("H94").Value = 2610.4443
$ws.Range("I94").Value = 3700
$ws.Range("J94").Value = 1738.8
$ws.Range("K94").Value = 3700
$ws.Range("L94").Value = 1738.8
$ws.Range("M94").Value = -3249
$ws.Range("N94").Value = -2640.8
# Row 99
$ws.Range("H99").Value = 14542.75
$ws.Range("I99").Value = 1890.3334
$ws.Range("J99").Value = 52500
$ws.Range("K99").Value = 1890.3334
$ws.Range("L99").Value = 52500
$ws.Range("M99").Value = -392.3334
$ws.Range("N99").Value = -55496
# Row 122
$ws.Range("H122").Value = 4317.25
$ws.Range("I122").Value = 2681
$ws.Range("J122").Value = 5953.5
$ws.Range("K122").Value = 8043
$ws.Range("L122").Value = 17860.5
$ws.Range("M122").Value = -5593
$ws.Range("N122").Value = -22760.5
# Row 126
$ws.Range("H126").Value = 14542.75
$ws.Range("I126").Value = 1890.3334
$ws.Range("J126").Value = 52500
$ws.Range("K126").Value = 5671.0002
$ws.Range("L126").Value = 157500
$ws.Range("M126").Value = -3201.0002
$ws.Range("N126").Value = -162440
# Row 136
$ws.Range("H136").Value = 2417.25
$ws.Range("I136").Value = 2212.8
$ws.Range("K136").Value = 6638.400000000001
$ws.Range("M136").Value = -4088.400000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 973591.7
$ws.Range("I107").Value = 5501.5
$ws.Range("J107").Value = 1134940.1
$ws.Range("K107").Value = 16504.5
$ws.Range("L107").Value = 3404820.3
$ws.Range("M107").Value = -14584.5
$ws.Range("N107").Value = -3408660.3
# Row 131
$ws.Range("H131").Value = 854.39
$ws.Range("J131").Value = 857.9697
$ws.Range("L131").Value = 2573.9091
$ws.Range("N131").Value = -12653.9091
# Row 136
$ws.Range("H136").Value = 1656.8518
$ws.Range("I136").Value = 1455
$ws.Range("K136").Value = 4365
$ws.Range("M136").Value = 735

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2678.4546
$ws.Range("I102").Value = 2739.3845
$ws.Range("J102").Value = 2590.4443
$ws.Range("K102").Value = 2739.3845
$ws.Range("L102").Value = 2590.4443
$ws.Range("M102").Value = -1117.3845
$ws.Range("N102").Value = -5834.4443
# Row 122
$ws.Range("H122").Value = 2442.0312
$ws.Range("I122").Value = 1899.2609
$ws.Range("J122").Value = 3829.111
$ws.Range("K122").Value = 5697.7827
$ws.Range("L122").Value = 11487.333
$ws.Range("M122").Value = -3247.7827
$ws.Range("N122").Value = -16387.333
